$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 2500
$ws.Range("J7").Value = 2500
$ws.Range("L7").Value = 2500
$ws.Range("N7").Value = -2724
$ws.Range("H10").Value = 32500
$ws.Range("J10").Value = 32500
$ws.Range("L10").Value = 32500
$ws.Range("N10").Value = -33086
$ws.Range("H14").Value = 2500
$ws.Range("J14").Value = 2500
$ws.Range("L14").Value = 2500
$ws.Range("N14").Value = -2882
$ws.Range("H21").Value = 200
$ws.Range("I21").Value = 200
$ws.Range("K21").Value = 200
$ws.Range("M21").Value = 268
$ws.Range("H23").Value = 200
$ws.Range("I23").Value = 200
$ws.Range("K23").Value = 200
$ws.Range("M23").Value = 34
$ws.Range("H39").Value = 2381.1
$ws.Range("I39").Value = 296.4
$ws.Range("J39").Value = 4465.8
$ws.Range("K39").Value = 889.1999999999999
$ws.Range("L39").Value = 13397.4
$ws.Range("M39").Value = -593.1999999999999
$ws.Range("N39").Value = -13989.4
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H62").Value = 50004000
$ws.Range("I62").Value = 50004000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 50004000
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 50004000
$ws.Range("I65").Value = 50004000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 250020000
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H101").Value = 2950.05
$ws.Range("I101").Value = 1243.0714
$ws.Range("K101").Value = 3729.2142
$ws.Range("M101").Value = -2107.2142
$ws.Range("H106").Value = 5352078.5
$ws.Range("I106").Value = 6065545.5
$ws.Range("J106").Value = 1075
$ws.Range("K106").Value = 6065545.5
$ws.Range("L106").Value = 1075
$ws.Range("M106").Value = -6064914.5
$ws.Range("N106").Value = -2337
$ws.Range("H116").Value = 10630.154
$ws.Range("I116").Value = 12736
$ws.Range("J116").Value = 7260.8
$ws.Range("K116").Value = 12736
$ws.Range("L116").Value = 7260.8
$ws.Range("M116").Value = -9294
$ws.Range("N116").Value = -14144.8
$ws.Range("H132").Value = 4063.8333
$ws.Range("I132").Value = 4039.7334
$ws.Range("J132").Value = 4425.3335
$ws.Range("K132").Value = 12119.2002
$ws.Range("L132").Value = 13276.0005
$ws.Range("M132").Value = -9589.200199999999
$ws.Range("N132").Value = -18336.0005
$ws.Range("H134").Value = 32149.615
$ws.Range("I134").Value = 22000
$ws.Range("J134").Value = 33995
$ws.Range("K134").Value = 22000
$ws.Range("L134").Value = 33995
$ws.Range("M134").Value = -16930
$ws.Range("N134").Value = -44135
$ws.Range("H136").Value = 54883.332
$ws.Range("J136").Value = 54883.332
$ws.Range("L136").Value = 54883.332
$ws.Range("N136").Value = -65083.332
$ws.Range("H137").Value = 2297.8333
$ws.Range("I137").Value = 1570.7778
$ws.Range("J137").Value = 2734.0667
$ws.Range("K137").Value = 4712.3334
$ws.Range("L137").Value = 8202.2001
$ws.Range("M137").Value = -2162.3334
$ws.Range("N137").Value = -13302.2001
$ws.Range("H138").Value = 3136.3733
$ws.Range("I138").Value = 2589.9375
$ws.Range("J138").Value = 3543.0232
$ws.Range("K138").Value = 7769.8125
$ws.Range("L138").Value = 10629.0696
$ws.Range("M138").Value = -2629.8125
$ws.Range("N138").Value = -20909.0696
$ws.Range("H141").Value = 3098.1956
$ws.Range("I141").Value = 2577.442
$ws.Range("J141").Value = 10562.333
$ws.Range("K141").Value = 7732.326
$ws.Range("L141").Value = 31686.999
$ws.Range("M141").Value = -2552.326
$ws.Range("N141").Value = -42046.999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1023.2308
$ws.Range("I2").Value = 499.5
$ws.Range("J2").Value = 1472.1428
$ws.Range("K2").Value = 499.5
$ws.Range("L2").Value = 1472.1428
$ws.Range("M2").Value = -386.5
$ws.Range("N2").Value = -1698.1428
$ws.Range("H3").Value = 7000
$ws.Range("J3").Value = 7000
$ws.Range("L3").Value = 7000
$ws.Range("N3").Value = -7230
$ws.Range("H32").Value = 2725.4658
$ws.Range("I32").Value = 2061.2354
$ws.Range("J32").Value = 11759
$ws.Range("K32").Value = 2061.2354
$ws.Range("L32").Value = 11759
$ws.Range("M32").Value = -1774.2354
$ws.Range("N32").Value = -12333
$ws.Range("H45").Value = 2044.05
$ws.Range("I45").Value = 1503.5454
$ws.Range("J45").Value = 2704.6667
$ws.Range("K45").Value = 1503.5454
$ws.Range("L45").Value = 2704.6667
$ws.Range("M45").Value = -1126.5454
$ws.Range("N45").Value = -3458.6667
$ws.Range("H50").Value = 7294.1
$ws.Range("I50").Value = 1887
$ws.Range("J50").Value = 10898.833
$ws.Range("K50").Value = 1887
$ws.Range("L50").Value = 10898.833
$ws.Range("M50").Value = -1173
$ws.Range("N50").Value = -12326.833
$ws.Range("H61").Value = 5065.9
$ws.Range("I61").Value = 5707.375
$ws.Range("J61").Value = 2500
$ws.Range("K61").Value = 5707.375
$ws.Range("L61").Value = 2500
$ws.Range("M61").Value = -5495.375
$ws.Range("N61").Value = -2924
$ws.Range("H74").Value = 2144.85
$ws.Range("I74").Value = 1459.5385
$ws.Range("K74").Value = 1459.5385
$ws.Range("M74").Value = -585.5385000000001
$ws.Range("H77").Value = 2144.85
$ws.Range("I77").Value = 1459.5385
$ws.Range("K77").Value = 7297.692500000001
$ws.Range("M77").Value = -2929.692500000001
$ws.Range("H97").Value = 1688.1333
$ws.Range("I97").Value = 876
$ws.Range("J97").Value = 3312.4
$ws.Range("K97").Value = 876
$ws.Range("L97").Value = 3312.4
$ws.Range("M97").Value = -380
$ws.Range("N97").Value = -4304.4
$ws.Range("H102").Value = 1382.8438
$ws.Range("I102").Value = 1266.1613
$ws.Range("K102").Value = 1266.1613
$ws.Range("M102").Value = 355.8387
$ws.Range("H116").Value = 1023.2308
$ws.Range("I116").Value = 499.5
$ws.Range("J116").Value = 1472.1428
$ws.Range("K116").Value = 499.5
$ws.Range("L116").Value = 1472.1428
$ws.Range("M116").Value = 1794.5
$ws.Range("N116").Value = -6060.1428
$ws.Range("H132").Value = 1938.5
$ws.Range("I132").Value = 1767.0566
$ws.Range("K132").Value = 5301.1698
$ws.Range("M132").Value = -2771.1698
$ws.Range("H133").Value = 47574.06
$ws.Range("J133").Value = 47574.06
$ws.Range("L133").Value = 47574.06
$ws.Range("N133").Value = -52634.06
$ws.Range("H136").Value = 5065.9
$ws.Range("I136").Value = 5707.375
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 17122.125
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -14572.125
$ws.Range("N136").Value = -12600
$ws.Range("H139").Value = 70000
$ws.Range("J139").Value = 70000
$ws.Range("L139").Value = 70000
$ws.Range("N139").Value = -80280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1023.2308
$ws.Range("I3").Value = 499.5
$ws.Range("J3").Value = 1472.1428
$ws.Range("K3").Value = 499.5
$ws.Range("L3").Value = 1472.1428
$ws.Range("M3").Value = -385.5
$ws.Range("N3").Value = -1700.1428
$ws.Range("H10").Value = 4699
$ws.Range("I10").Value = 4699
$ws.Range("K10").Value = 4699
$ws.Range("M10").Value = -4559
$ws.Range("H20").Value = 2319
$ws.Range("I20").Value = 1149.5
$ws.Range("J20").Value = 6997
$ws.Range("K20").Value = 1149.5
$ws.Range("L20").Value = 6997
$ws.Range("M20").Value = -902.5
$ws.Range("N20").Value = -7491
$ws.Range("H35").Value = 79999
$ws.Range("J35").Value = 79999
$ws.Range("L35").Value = 79999
$ws.Range("N35").Value = -80619
$ws.Range("H94").Value = 2383902
$ws.Range("I94").Value = 2703282
$ws.Range("K94").Value = 2703282
$ws.Range("M94").Value = -2702831
$ws.Range("H105").Value = 3758.7693
$ws.Range("I105").Value = 3533.0908
$ws.Range("K105").Value = 3533.0908
$ws.Range("M105").Value = -1786.0908
$ws.Range("H107").Value = 5160.643
$ws.Range("I107").Value = 5386.6924
$ws.Range("K107").Value = 5386.6924
$ws.Range("M107").Value = -3466.6924
$ws.Range("H133").Value = 85333.336
$ws.Range("I133").Value = 60000
$ws.Range("J133").Value = 98000
$ws.Range("K133").Value = 60000
$ws.Range("L133").Value = 98000
$ws.Range("M133").Value = -54940
$ws.Range("N133").Value = -108120
$ws.Range("H134").Value = 5077.0327
$ws.Range("I134").Value = 4391.75
$ws.Range("J134").Value = 6850.706
$ws.Range("K134").Value = 13175.25
$ws.Range("L134").Value = 20552.118
$ws.Range("M134").Value = -10640.25
$ws.Range("N134").Value = -25622.118

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 917.2
$ws.Range("I22").Value = 171.5
$ws.Range("K22").Value = 171.5
$ws.Range("M22").Value = 178.5
$ws.Range("H29").Value = 19499
$ws.Range("J29").Value = 23332
$ws.Range("L29").Value = 23332
$ws.Range("N29").Value = -23918
$ws.Range("H31").Value = 4594.1963
$ws.Range("I31").Value = 2611
$ws.Range("J31").Value = 5204.41
$ws.Range("K31").Value = 2611
$ws.Range("L31").Value = 5204.41
$ws.Range("M31").Value = -2316
$ws.Range("N31").Value = -5794.41
$ws.Range("H34").Value = 4594.1963
$ws.Range("I34").Value = 2611
$ws.Range("J34").Value = 5204.41
$ws.Range("K34").Value = 2611
$ws.Range("L34").Value = 5204.41
$ws.Range("M34").Value = -2409
$ws.Range("N34").Value = -5608.41
$ws.Range("H94").Value = 2773.375
$ws.Range("I94").Value = 909.5
$ws.Range("K94").Value = 909.5
$ws.Range("M94").Value = -458.5
$ws.Range("H99").Value = 1896.875
$ws.Range("I99").Value = 1882.1428
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 1882.1428
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -384.1428000000001
$ws.Range("N99").Value = -4996
$ws.Range("H105").Value = 819.6667
$ws.Range("I105").Value = 825.75
$ws.Range("K105").Value = 825.75
$ws.Range("M105").Value = 921.25
$ws.Range("H126").Value = 1896.875
$ws.Range("I126").Value = 1882.1428
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 5646.428400000001
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -3176.428400000001
$ws.Range("N126").Value = -10940
$ws.Range("H132").Value = 3935.4546
$ws.Range("I132").Value = 3673.3142
$ws.Range("J132").Value = 4954.8887
$ws.Range("K132").Value = 11019.9426
$ws.Range("L132").Value = 14864.6661
$ws.Range("M132").Value = -8489.942599999998
$ws.Range("N132").Value = -19924.6661
$ws.Range("H134").Value = 3490.6943
$ws.Range("I134").Value = 2677.3125
$ws.Range("J134").Value = 9997.75
$ws.Range("K134").Value = 8031.9375
$ws.Range("L134").Value = 29993.25
$ws.Range("M134").Value = -5496.9375
$ws.Range("N134").Value = -35063.25
$ws.Range("H141").Value = 60000
$ws.Range("J141").Value = 60000
$ws.Range("L141").Value = 60000
$ws.Range("N141").Value = -70360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2979.9092
$ws.Range("I3").Value = 1442.1111
$ws.Range("K3").Value = 4326.3333
$ws.Range("M3").Value = -4214.3333
$ws.Range("H32").Value = 376249.5
$ws.Range("I32").Value = 334666.5
$ws.Range("J32").Value = 500998.5
$ws.Range("K32").Value = 1003999.5
$ws.Range("L32").Value = 1502995.5
$ws.Range("M32").Value = -1003716.5
$ws.Range("N32").Value = -1503561.5
$ws.Range("H60").Value = 736.4286
$ws.Range("I60").Value = 754.61536
$ws.Range("J60").Value = 500
$ws.Range("K60").Value = 2263.84608
$ws.Range("L60").Value = 1500
$ws.Range("M60").Value = -2012.84608
$ws.Range("N60").Value = -2002
$ws.Range("H70").Value = 16666.666
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 16666.666
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -50629.99800000001
$ws.Range("H73").Value = 16666.666
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 16666.666
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -52183.99800000001
$ws.Range("H76").Value = 13285.571
$ws.Range("I76").Value = 10499.833
$ws.Range("J76").Value = 30000
$ws.Range("K76").Value = 31499.499
$ws.Range("L76").Value = 90000
$ws.Range("M76").Value = -31116.499
$ws.Range("N76").Value = -90766
$ws.Range("H79").Value = 13285.571
$ws.Range("I79").Value = 10499.833
$ws.Range("J79").Value = 30000
$ws.Range("K79").Value = 31499.499
$ws.Range("L79").Value = 90000
$ws.Range("M79").Value = -30173.499
$ws.Range("N79").Value = -92652
$ws.Range("H97").Value = 782.875
$ws.Range("I97").Value = 316
$ws.Range("J97").Value = 1249.75
$ws.Range("K97").Value = 948
$ws.Range("L97").Value = 3749.25
$ws.Range("M97").Value = -452
$ws.Range("N97").Value = -4741.25
$ws.Range("H121").Value = 17615446
$ws.Range("I121").Value = 34366.668
$ws.Range("J121").Value = 20911898
$ws.Range("K121").Value = 103100.004
$ws.Range("L121").Value = 62735694
$ws.Range("M121").Value = -101790.004
$ws.Range("N121").Value = -62738314
$ws.Range("H123").Value = 14517.182
$ws.Range("I123").Value = 1038
$ws.Range("J123").Value = 25749.834
$ws.Range("K123").Value = 3114
$ws.Range("L123").Value = 77249.50199999999
$ws.Range("M123").Value = -664
$ws.Range("N123").Value = -82149.50199999999
$ws.Range("H134").Value = 14892.667
$ws.Range("I134").Value = 5678
$ws.Range("J134").Value = 19500
$ws.Range("K134").Value = 17034
$ws.Range("L134").Value = 58500
$ws.Range("M134").Value = -11964
$ws.Range("N134").Value = -68640
$ws.Range("H139").Value = 10424451
$ws.Range("I139").Value = 11909373
$ws.Range("K139").Value = 35728119
$ws.Range("M139").Value = -35722979
$ws.Range("H140").Value = 5565319.5
$ws.Range("I140").Value = 25003226
$ws.Range("J140").Value = 11631.6
$ws.Range("K140").Value = 75009678
$ws.Range("L140").Value = 34894.8
$ws.Range("M140").Value = -75004498
$ws.Range("N140").Value = -45254.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 17709.572
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 17709.572
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -18213.572
$ws.Range("H52").Value = 999999.5
$ws.Range("I52").Value = 999999.5
$ws.Range("K52").Value = 999999.5
$ws.Range("M52").Value = -999740.5
$ws.Range("H70").Value = 7608.6665
$ws.Range("I70").Value = 7323.385
$ws.Range("J70").Value = 8350.4
$ws.Range("K70").Value = 7323.385
$ws.Range("L70").Value = 8350.4
$ws.Range("M70").Value = -7053.385
$ws.Range("N70").Value = -8890.4
$ws.Range("H73").Value = 7608.6665
$ws.Range("I73").Value = 7323.385
$ws.Range("J73").Value = 8350.4
$ws.Range("K73").Value = 7323.385
$ws.Range("L73").Value = 8350.4
$ws.Range("M73").Value = -6387.385
$ws.Range("N73").Value = -10222.4
$ws.Range("H80").Value = 2088.9167
$ws.Range("I80").Value = 2353.4285
$ws.Range("J80").Value = 1718.6
$ws.Range("K80").Value = 2353.4285
$ws.Range("L80").Value = 1718.6
$ws.Range("M80").Value = -1355.4285
$ws.Range("N80").Value = -3714.6
$ws.Range("H83").Value = 2088.9167
$ws.Range("I83").Value = 2353.4285
$ws.Range("J83").Value = 1718.6
$ws.Range("K83").Value = 11767.1425
$ws.Range("L83").Value = 8593
$ws.Range("M83").Value = -6775.1425
$ws.Range("N83").Value = -18577
$ws.Range("H97").Value = 1200.9524
$ws.Range("I97").Value = 536
$ws.Range("J97").Value = 2281.5
$ws.Range("K97").Value = 536
$ws.Range("L97").Value = 2281.5
$ws.Range("M97").Value = -40
$ws.Range("N97").Value = -3273.5
$ws.Range("H102").Value = 3271.9167
$ws.Range("I102").Value = 2980.4092
$ws.Range("J102").Value = 3730
$ws.Range("K102").Value = 2980.4092
$ws.Range("L102").Value = 3730
$ws.Range("M102").Value = -1358.4092
$ws.Range("N102").Value = -6974
$ws.Range("H126").Value = 6799.091
$ws.Range("J126").Value = 5750
$ws.Range("L126").Value = 17250
$ws.Range("N126").Value = -22190
$ws.Range("H132").Value = 1864.4222
$ws.Range("I132").Value = 1651.6757
$ws.Range("J132").Value = 2848.375
$ws.Range("K132").Value = 4955.0271
$ws.Range("L132").Value = 8545.125
$ws.Range("M132").Value = -2425.0271
$ws.Range("N132").Value = -13605.125
$ws.Range("H134").Value = 21662
$ws.Range("J134").Value = 21662
$ws.Range("L134").Value = 64986
$ws.Range("N134").Value = -70056

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 104899.1
$ws.Range("I7").Value = 147320.14
$ws.Range("J7").Value = 5916.6665
$ws.Range("K7").Value = 147320.14
$ws.Range("L7").Value = 5916.6665
$ws.Range("M7").Value = -147208.14
$ws.Range("N7").Value = -6140.6665
$ws.Range("H22").Value = 1109.7693
$ws.Range("I22").Value = 841.5
$ws.Range("K22").Value = 841.5
$ws.Range("M22").Value = -546.5
$ws.Range("H27").Value = 1109.7693
$ws.Range("I27").Value = 841.5
$ws.Range("K27").Value = 841.5
$ws.Range("M27").Value = -734.5
$ws.Range("H55").Value = 420.075
$ws.Range("I55").Value = 342.66666
$ws.Range("J55").Value = 652.3
$ws.Range("K55").Value = 342.66666
$ws.Range("L55").Value = 652.3
$ws.Range("M55").Value = -169.66666
$ws.Range("N55").Value = -998.3
$ws.Range("H61").Value = 2500.2856
$ws.Range("I61").Value = 2500.2856
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2500.2856
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H68").Value = 2284.3809
$ws.Range("I68").Value = 2305.6155
$ws.Range("J68").Value = 2249.875
$ws.Range("K68").Value = 2305.6155
$ws.Range("L68").Value = 2249.875
$ws.Range("M68").Value = -1556.6155
$ws.Range("N68").Value = -3747.875
$ws.Range("H71").Value = 2284.3809
$ws.Range("I71").Value = 2305.6155
$ws.Range("J71").Value = 2249.875
$ws.Range("K71").Value = 11528.0775
$ws.Range("L71").Value = 11249.375
$ws.Range("M71").Value = -7784.077499999999
$ws.Range("N71").Value = -18737.375
$ws.Range("H82").Value = 2767.818
$ws.Range("I82").Value = 2281.125
$ws.Range("J82").Value = 4065.6667
$ws.Range("K82").Value = 2281.125
$ws.Range("L82").Value = 4065.6667
$ws.Range("M82").Value = -1920.125
$ws.Range("N82").Value = -4787.6667
$ws.Range("H85").Value = 2767.818
$ws.Range("I85").Value = 2281.125
$ws.Range("J85").Value = 4065.6667
$ws.Range("K85").Value = 2281.125
$ws.Range("L85").Value = 4065.6667
$ws.Range("M85").Value = -1033.125
$ws.Range("N85").Value = -6561.6667
$ws.Range("I93").Value = 1045.8182
$ws.Range("J93").Value = 2166.3333
$ws.Range("K93").Value = 1045.8182
$ws.Range("L93").Value = 2166.3333
$ws.Range("M93").Value = 202.1818000000001
$ws.Range("N93").Value = -4662.3333
$ws.Range("H113").Value = 2500.2856
$ws.Range("I113").Value = 2500.2856
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2500.2856
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H126").Value = 104899.1
$ws.Range("I126").Value = 147320.14
$ws.Range("J126").Value = 5916.6665
$ws.Range("K126").Value = 441960.42
$ws.Range("L126").Value = 17749.9995
$ws.Range("M126").Value = -439490.42
$ws.Range("N126").Value = -22689.9995
$ws.Range("H132").Value = 1986.2106
$ws.Range("I132").Value = 1679.5834
$ws.Range("J132").Value = 2511.8572
$ws.Range("K132").Value = 5038.7502
$ws.Range("L132").Value = 7535.571599999999
$ws.Range("M132").Value = -2508.7502
$ws.Range("N132").Value = -12595.5716
$ws.Range("H136").Value = 1270.0807
$ws.Range("I136").Value = 1118.46
$ws.Range("K136").Value = 3355.38
$ws.Range("M136").Value = -805.3800000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 5500
$ws.Range("J7").Value = 10000
$ws.Range("L7").Value = 10000
$ws.Range("N7").Value = -10226
$ws.Range("H9").Value = 2226.5
$ws.Range("J9").Value = 3500
$ws.Range("L9").Value = 3500
$ws.Range("N9").Value = -3780
$ws.Range("H14").Value = 17927
$ws.Range("I14").Value = 20164.834
$ws.Range("J14").Value = 4500
$ws.Range("K14").Value = 20164.834
$ws.Range("L14").Value = 4500
$ws.Range("M14").Value = -19996.834
$ws.Range("N14").Value = -4836
$ws.Range("H46").Value = 49388.125
$ws.Range("J46").Value = 49388.125
$ws.Range("L46").Value = 49388.125
$ws.Range("N46").Value = -49850.125
$ws.Range("H96").Value = 4999.5
$ws.Range("J96").Value = 4999.5
$ws.Range("L96").Value = 4999.5
$ws.Range("N96").Value = -7745.5
$ws.Range("H100").Value = 864
$ws.Range("I100").Value = 864
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1728
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H132").Value = 2065.7097
$ws.Range("I132").Value = 2063.9873
$ws.Range("J132").Value = 2075.4285
$ws.Range("K132").Value = 6191.961899999999
$ws.Range("L132").Value = 6226.2855
$ws.Range("M132").Value = -3661.961899999999
$ws.Range("N132").Value = -11286.2855
$ws.Range("H134").Value = 49388.125
$ws.Range("J134").Value = 49388.125
$ws.Range("L134").Value = 148164.375
$ws.Range("N134").Value = -153234.375
$ws.Range("H136").Value = 2735.2693
$ws.Range("I136").Value = 3439.25
$ws.Range("J136").Value = 1608.9
$ws.Range("K136").Value = 10317.75
$ws.Range("L136").Value = 4826.700000000001
$ws.Range("M136").Value = -7767.75
$ws.Range("N136").Value = -9926.700000000001
